$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so numeric-looking values
# (e.g. "351.16", "0.0834", "1.00") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.526.67"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.780.46"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "351.16"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "108.40"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +4.78%  "
$ws.Range("D10").Value = "39.30"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "0.0834"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "19.84"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").Value = "3.215.07"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "2.767.02"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "0.926"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "51.519.60"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +4.53%  "
$ws.Range("D20").Value = "3.09"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "70.31"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "266.62"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "25.90"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("D28").Value = "0.165"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "10.28"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "37.17"
$ws.Range("E30").Value = "  +10.36%  "
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "6.18"
$ws.Range("E32").Value = "  +7.99%  "
$ws.Range("D33").Value = "51.95"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "0.0447"
$ws.Range("E34").Value = "  -5.31%  "
$ws.Range("D35").Value = "5.56"
$ws.Range("E35").Value = "  +7.00%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "0.0837"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "18.66"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("D43").Value = "120.39"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "2.19"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").Value = "2.148.96"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("D49").Value = "0.224"
$ws.Range("E49").Value = "  +17.89%  "
$ws.Range("D50").Value = "5.48"
$ws.Range("E50").Value = "  -5.91%  "
$ws.Range("D51").Value = "0.904"
$ws.Range("E51").Value = "  -5.17%  "
